$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typos in the two column headers (FR / NFR) ---
$ws.Range("B3").Value = "requirements of place order - FR"
$ws.Range("I3").Value = "requirements of place order - NFR"

# --- Reword several existing FR (functional requirement) rows in column C ---
$ws.Range("C16").Value = "the system will allow a customer to enter the waiting list "
$ws.Range("C17").Value = "the system will be able to recognize how much time has been passed"
$ws.Range("C18").Value = "the system will be able to generate codes"
$ws.Range("C19").Value = "the ststem will be able to produce bills "
$ws.Range("C20").Value = "the system will be able to distinguish between customers types"
$ws.Range("C21").Value = "the system will enable more options for a subscriber  "
$ws.Range("C22").Value = "the system will be able to change a table status"

# --- Remove the 20th FR row (B23/C23): its text becomes the new NFR #20 entry ---
$ws.Range("B23").Clear()
$ws.Range("C23").ClearContents()

# --- NFR #20 (row 23) now holds the text that used to be the removed FR item's neighbor ---
$ws.Range("J23").Value = "all customers data will be stored securely"

# --- Append 4 new NFR rows (21-24) split out from former combined NFR entries ---
# Copy formatting (style) from the last existing NFR row down into the new rows first.
$ws.Range("I23:J23").Copy($ws.Range("I24:J24"))
$ws.Range("I23:J23").Copy($ws.Range("I25:J25"))
$ws.Range("I23:J23").Copy($ws.Range("I26:J26"))
$ws.Range("I23:J23").Copy($ws.Range("I27:J27"))

$ws.Range("I24").Value = 21
$ws.Range("J24").Value = "entry to the waiting list depends if there isnt a free table"

$ws.Range("I25").Value = 22
$ws.Range("J25").Value = "if the cutomer lost his code the system will regenerate a new one"

$ws.Range("I26").Value = 23
$ws.Range("J26").Value = "subscriber options are: view his order history/visits.. "

$ws.Range("I27").Value = 24
$ws.Range("J27").Value = "table staus can be: occupied/available"
